$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: R15 footprint/comment/part# update
$ws.Range("B30").Value = "R_0603_1608Metric"
$ws.Range("C30").Value = "47.5K"
$ws.Range("D30").Value = "C137715"

# Row 34: R14 footprint update (comment & part# unchanged)
$ws.Range("B34").Value = "R_0603_1608Metric"

# Move the active selection to D30 (matches the author's saved cursor position)
[void]$ws.Range("D30").Select()
